# RPA datasets push 2024-07-05
#
# The "뱅크웨어글로벌" (Bankware Global) IPO row, previously listed between
# "산일전기(유가)" and "이베스트스팩6호" with a demand-forecast date of
# 2024.07.08~07.12, is moved up so it now follows "키움스팩9호" (its demand
# forecast date also updates to 2024.07.23~07.29). Every row that used to sit
# between the old and new positions shifts down by one row to close the gap.
#
# Net effect on row contents (rows 14-21 are untouched):
#   Row 5  -> 뱅크웨어글로벌 (moved up, date changed to 2024.07.23~07.29)
#   Row 6  -> 티디에스팜 (was row 5)
#   Row 7  -> 아이빔테크놀로지 (was row 6)
#   Row 8  -> 넥스트바이오메디컬 (was row 7)
#   Row 9  -> 피앤에스미캐닉스 (was row 8)
#   Row 10 -> 케이쓰리아이 (was row 9)
#   Row 11 -> NH스팩31호 (was row 10)
#   Row 12 -> SK증권스팩13호 (was row 11)
#   Row 13 -> 산일전기(유가) (was row 12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 5;  A = "뱅크웨어글로벌";      B = "2024.07.23~07.29"; C = "16,000~19,000"; D = "-"; E = 22400;  F = "미래에셋증권" },
    @{ Row = 6;  A = "티디에스팜";          B = "2024.07.18~07.24"; C = "9,500~10,700";  D = "-"; E = 9500;   F = "한국투자증권" },
    @{ Row = 7;  A = "아이빔테크놀로지";     B = "2024.07.15~07.19"; C = "7,300~8,500";   D = "-"; E = 16308;  F = "삼성증권" },
    @{ Row = 8;  A = "넥스트바이오메디컬";   B = "2024.07.15~07.19"; C = "24,000~29,000"; D = "-"; E = 24000;  F = "한국투자증권" },
    @{ Row = 9;  A = "피앤에스미캐닉스";     B = "2024.07.11~07.17"; C = "14,000~17,000"; D = "-"; E = 18900;  F = "키움증권" },
    @{ Row = 10; A = "케이쓰리아이";        B = "2024.07.10~07.16"; C = "12,500~15,500"; D = "-"; E = 17500;  F = "하나증권" },
    @{ Row = 11; A = "NH스팩31호";         B = "2024.07.09~07.10"; C = "2,000~2,000";   D = "-"; E = 12000;  F = "NH투자증권" },
    @{ Row = 12; A = "SK증권스팩13호";      B = "2024.07.09~07.10"; C = "2,000~2,000";   D = "-"; E = 8000;   F = "SK증권" },
    @{ Row = 13; A = "산일전기(유가)";      B = "2024.07.09~07.15"; C = "24,000~30,000"; D = "-"; E = 182400; F = "미래에셋증권,삼성증권" }
)

foreach ($r in $rows) {
    $ws.Range("A$($r.Row)").Value = $r.A
    $ws.Range("B$($r.Row)").Value = $r.B
    $ws.Range("C$($r.Row)").Value = $r.C
    $ws.Range("D$($r.Row)").Value = $r.D
    $ws.Range("E$($r.Row)").Value = $r.E
    $ws.Range("F$($r.Row)").Value = $r.F
}
